# Fix financial data rows 2-9 (columns D:AJ) of the IFRS company list sheet.
# Each row is rewritten in one shot via a 2D Variant array assigned to the
# Range.Value; $null entries clear the corresponding cell (matching cells
# that are removed entirely from the sheet, e.g. J2/O2/Y2/Z2/AD2/AH2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$arr2 = New-Object "object[,]" 1,33
$arr2[0,0] = 1652
$arr2[0,1] = 99
$arr2[0,2] = 99
$arr2[0,3] = 220
$arr2[0,4] = 202
$arr2[0,5] = 202
$arr2[0,6] = $null
$arr2[0,7] = 2164
$arr2[0,8] = 1128
$arr2[0,9] = 1036
$arr2[0,10] = 1036
$arr2[0,11] = $null
$arr2[0,12] = 750
$arr2[0,13] = 104
$arr2[0,14] = -1312
$arr2[0,15] = 1441
$arr2[0,16] = 108
$arr2[0,17] = -4
$arr2[0,18] = 608
$arr2[0,19] = 6.01
$arr2[0,20] = 12.25
$arr2[0,21] = $null
$arr2[0,22] = $null
$arr2[0,23] = 108.92
$arr2[0,24] = 38.1
$arr2[0,25] = 1483
$arr2[0,26] = $null
$arr2[0,27] = 6905
$arr2[0,28] = 0
$arr2[0,29] = 0
$arr2[0,30] = $null
$arr2[0,31] = 0
$arr2[0,32] = 15000000
$ws.Range("D2:AJ2").Value = $arr2

# Row 3
$arr3 = New-Object "object[,]" 1,33
$arr3[0,0] = 2460
$arr3[0,1] = 188
$arr3[0,2] = 188
$arr3[0,3] = 178
$arr3[0,4] = 147
$arr3[0,5] = 147
$arr3[0,6] = $null
$arr3[0,7] = 2292
$arr3[0,8] = 1102
$arr3[0,9] = 1189
$arr3[0,10] = 1189
$arr3[0,11] = $null
$arr3[0,12] = 750
$arr3[0,13] = 298
$arr3[0,14] = -132
$arr3[0,15] = -1
$arr3[0,16] = 118
$arr3[0,17] = 180
$arr3[0,18] = 606
$arr3[0,19] = 7.64
$arr3[0,20] = 5.97
$arr3[0,21] = 13.2
$arr3[0,22] = 6.59
$arr3[0,23] = 92.69
$arr3[0,24] = 58.57
$arr3[0,25] = 979
$arr3[0,26] = $null
$arr3[0,27] = 7929
$arr3[0,28] = 0
$arr3[0,29] = 0
$arr3[0,30] = $null
$arr3[0,31] = 0
$arr3[0,32] = 15000000
$ws.Range("D3:AJ3").Value = $arr3

# Row 4
$arr4 = New-Object "object[,]" 1,33
$arr4[0,0] = 2762
$arr4[0,1] = 258
$arr4[0,2] = 258
$arr4[0,3] = 232
$arr4[0,4] = 188
$arr4[0,5] = 189
$arr4[0,6] = 0
$arr4[0,7] = 2510
$arr4[0,8] = 896
$arr4[0,9] = 1614
$arr4[0,10] = 1611
$arr4[0,11] = 2
$arr4[0,12] = 850
$arr4[0,13] = 164
$arr4[0,14] = -621
$arr4[0,15] = 212
$arr4[0,16] = 307
$arr4[0,17] = -142
$arr4[0,18] = 592
$arr4[0,19] = 9.35
$arr4[0,20] = 6.83
$arr4[0,21] = 13.47
$arr4[0,22] = 7.85
$arr4[0,23] = 55.53
$arr4[0,24] = 89.55
$arr4[0,25] = 1174
$arr4[0,26] = 11.03
$arr4[0,27] = 9477
$arr4[0,28] = 1.37
$arr4[0,29] = 250
$arr4[0,30] = 1.93
$arr4[0,31] = 22.53
$arr4[0,32] = 17000000
$ws.Range("D4:AJ4").Value = $arr4

# Row 5
$arr5 = New-Object "object[,]" 1,33
$arr5[0,0] = 3251
$arr5[0,1] = 339
$arr5[0,2] = 339
$arr5[0,3] = 294
$arr5[0,4] = 239
$arr5[0,5] = 240
$arr5[0,6] = -1
$arr5[0,7] = 2904
$arr5[0,8] = 1093
$arr5[0,9] = 1811
$arr5[0,10] = 1809
$arr5[0,11] = 1
$arr5[0,12] = 850
$arr5[0,13] = 237
$arr5[0,14] = -365
$arr5[0,15] = 106
$arr5[0,16] = 365
$arr5[0,17] = -128
$arr5[0,18] = 726
$arr5[0,19] = 10.42
$arr5[0,20] = 7.35
$arr5[0,21] = 14.02
$arr5[0,22] = 8.83
$arr5[0,23] = 60.37
$arr5[0,24] = 112.88
$arr5[0,25] = 1410
$arr5[0,26] = 11.7
$arr5[0,27] = 10644
$arr5[0,28] = 1.55
$arr5[0,29] = 300
$arr5[0,30] = 1.82
$arr5[0,31] = 21.27
$arr5[0,32] = 17000000
$ws.Range("D5:AJ5").Value = $arr5

# Row 6
$arr6 = New-Object "object[,]" 1,33
$arr6[0,0] = 3633
$arr6[0,1] = 268
$arr6[0,2] = 268
$arr6[0,3] = 245
$arr6[0,4] = 207
$arr6[0,5] = 208
$arr6[0,6] = $null
$arr6[0,7] = 3062
$arr6[0,8] = 1099
$arr6[0,9] = 1964
$arr6[0,10] = 1963
$arr6[0,11] = $null
$arr6[0,12] = 850
$arr6[0,13] = 368
$arr6[0,14] = -296
$arr6[0,15] = -105
$arr6[0,16] = 299
$arr6[0,17] = 69
$arr6[0,18] = 678
$arr6[0,19] = 7.37
$arr6[0,20] = 5.71
$arr6[0,21] = 11.05
$arr6[0,22] = 6.95
$arr6[0,23] = 55.96
$arr6[0,24] = 130.98
$arr6[0,25] = 1226
$arr6[0,26] = 10.56
$arr6[0,27] = 11549
$arr6[0,28] = 1.12
$arr6[0,29] = $null
$arr6[0,30] = $null
$arr6[0,31] = 24.47
$arr6[0,32] = 17000000
$ws.Range("D6:AJ6").Value = $arr6

# Row 7
$arr7 = New-Object "object[,]" 1,33
$arr7[0,0] = 3772
$arr7[0,1] = 240
$arr7[0,2] = $null
$arr7[0,3] = 229
$arr7[0,4] = 193
$arr7[0,5] = 194
$arr7[0,6] = $null
$arr7[0,7] = 3350
$arr7[0,8] = 1245
$arr7[0,9] = 2105
$arr7[0,10] = 2105
$arr7[0,11] = $null
$arr7[0,12] = 850
$arr7[0,13] = 374
$arr7[0,14] = -233
$arr7[0,15] = 5
$arr7[0,16] = 219
$arr7[0,17] = 160
$arr7[0,18] = $null
$arr7[0,19] = 6.36
$arr7[0,20] = 5.12
$arr7[0,21] = 9.52
$arr7[0,22] = 6.02
$arr7[0,23] = 59.17
$arr7[0,24] = $null
$arr7[0,25] = 1139
$arr7[0,26] = 16.2
$arr7[0,27] = 12384
$arr7[0,28] = 1.49
$arr7[0,29] = 317
$arr7[0,30] = 1.72
$arr7[0,31] = 27.8
$arr7[0,32] = $null
$ws.Range("D7:AJ7").Value = $arr7

# Row 8
$arr8 = New-Object "object[,]" 1,33
$arr8[0,0] = 4336
$arr8[0,1] = 405
$arr8[0,2] = $null
$arr8[0,3] = 375
$arr8[0,4] = 314
$arr8[0,5] = 314
$arr8[0,6] = $null
$arr8[0,7] = 3650
$arr8[0,8] = 1291
$arr8[0,9] = 2359
$arr8[0,10] = 2361
$arr8[0,11] = $null
$arr8[0,12] = 850
$arr8[0,13] = 469
$arr8[0,14] = -198
$arr8[0,15] = -56
$arr8[0,16] = 193
$arr8[0,17] = 268
$arr8[0,18] = $null
$arr8[0,19] = 9.33
$arr8[0,20] = 7.24
$arr8[0,21] = 14.04
$arr8[0,22] = 8.960000000000001
$arr8[0,23] = 54.7
$arr8[0,24] = $null
$arr8[0,25] = 1847
$arr8[0,26] = 9.66
$arr8[0,27] = 13888
$arr8[0,28] = 1.29
$arr8[0,29] = 393
$arr8[0,30] = 2.2
$arr8[0,31] = 21.3
$arr8[0,32] = $null
$ws.Range("D8:AJ8").Value = $arr8

# Row 9
$arr9 = New-Object "object[,]" 1,33
$arr9[0,0] = 4813
$arr9[0,1] = 454
$arr9[0,2] = $null
$arr9[0,3] = 421
$arr9[0,4] = 351
$arr9[0,5] = 351
$arr9[0,6] = $null
$arr9[0,7] = 4063
$arr9[0,8] = 1407
$arr9[0,9] = 2656
$arr9[0,10] = 2655
$arr9[0,11] = $null
$arr9[0,12] = 850
$arr9[0,13] = 463
$arr9[0,14] = -218
$arr9[0,15] = -59
$arr9[0,16] = 211
$arr9[0,17] = 244
$arr9[0,18] = $null
$arr9[0,19] = 9.43
$arr9[0,20] = 7.29
$arr9[0,21] = 13.99
$arr9[0,22] = 9.09
$arr9[0,23] = 52.96
$arr9[0,24] = $null
$arr9[0,25] = 2065
$arr9[0,26] = 8.65
$arr9[0,27] = 15620
$arr9[0,28] = 1.14
$arr9[0,29] = 423
$arr9[0,30] = 2.37
$arr9[0,31] = 20.5
$arr9[0,32] = $null
$ws.Range("D9:AJ9").Value = $arr9
